$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the data block from A4:B5 up to A1:B2
$src = $ws.Range("A4:B5")
$dst = $ws.Range("A1:B2")
$dst.Value2 = $src.Value2

# Clear out the old rows now that the data has moved
$src.ClearContents()

# Left-align the numeric id cell that moved to A2
$ws.Range("A2").HorizontalAlignment = -4131

# Update selection / active cell to match the saved view state
$ws.Range("C4").Select()
